$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table: 4 sending clusters (ECs, FAPs, M2, sCs) x 4 target clusters (ECs, FAPs, M2, sCs)
$data = @(
  @("ECs", "Fgf2", "Sdc4", "ECs", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 3.001642333333333, 9.004927, 0.05169795991651582, 0.05169795991651582, 0.645876387980111, 5.812887491821, 0.001005984779176297, 0.001005984779176297),
  @("ECs", "Fgf2", "Sdc4", "FAPs", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 11.05428166666667, 33.162845, 0.1903903753498087, 0.1903903753498088, 2.378597688103889, 21.407379192935, 0.003704784869903195, 0.003704784869903196),
  @("ECs", "Fgf2", "Sdc4", "M2", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 16.49405866666666, 49.482176, 0.2840808761059341, 0.2840808761059341, 3.549098077560888, 31.94188269804799, 0.00552789777157801, 0.005527897771578011),
  @("ECs", "Fgf2", "Sdc4", "sCs", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 27.51115433333333, 82.533463, 0.4738307886277414, 0.4738307886277414, 5.919694292905445, 53.277248636149, 0.009220219947447667, 0.009220219947447667),
  @("FAPs", "Fgf2", "Sdc4", "ECs", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 3.001642333333333, 9.004927, 0.05169795991651582, 0.05169795991651582, 24.24480942414311, 218.203284817288, 0.03776250333441369, 0.03776250333441369),
  @("FAPs", "Fgf2", "Sdc4", "FAPs", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 11.05428166666667, 33.162845, 0.1903903753498087, 0.1903903753498088, 89.28743753140891, 803.5869377826801, 0.1390696498584769, 0.1390696498584769),
  @("FAPs", "Fgf2", "Sdc4", "M2", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 16.49405866666666, 49.482176, 0.2840808761059341, 0.2840808761059341, 133.2255027733049, 1199.029524959744, 0.2075053841296043, 0.2075053841296043),
  @("FAPs", "Fgf2", "Sdc4", "sCs", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 27.51115433333333, 82.533463, 0.4738307886277414, 0.4738307886277414, 222.2125822396524, 1999.913240156872, 0.3461072112786932, 0.3461072112786933),
  @("M2", "Fgf2", "Sdc4", "ECs", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 3.001642333333333, 9.004927, 0.05169795991651582, 0.05169795991651582, 1.558380660050667, 14.025425940456, 0.002427255823790786, 0.002427255823790787),
  @("M2", "Fgf2", "Sdc4", "FAPs", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 11.05428166666667, 33.162845, 0.1903903753498087, 0.1903903753498088, 5.739117738573333, 51.65205964716001, 0.008938962932150494, 0.008938962932150497),
  @("M2", "Fgf2", "Sdc4", "M2", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 16.49405866666666, 49.482176, 0.2840808761059341, 0.2840808761059341, 8.563319402325332, 77.06987462092799, 0.01333779828196727, 0.01333779828196727),
  @("M2", "Fgf2", "Sdc4", "sCs", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 27.51115433333333, 82.533463, 0.4738307886277414, 0.4738307886277414, 14.28313106216267, 128.548179559464, 0.02224669103085137, 0.02224669103085138),
  @("sCs", "Fgf2", "Sdc4", "ECs", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 3.001642333333333, 9.004927, 0.05169795991651582, 0.05169795991651582, 6.742779277381111, 60.68501349643, 0.01050221597913505, 0.01050221597913505),
  @("sCs", "Fgf2", "Sdc4", "FAPs", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 11.05428166666667, 33.162845, 0.1903903753498087, 0.1903903753498088, 24.83193301233889, 223.48739711105, 0.03867697768927819, 0.0386769776892782),
  @("sCs", "Fgf2", "Sdc4", "M2", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 16.49405866666666, 49.482176, 0.2840808761059341, 0.2840808761059341, 37.05164860664888, 333.46483745984, 0.05770979592278456, 0.05770979592278457),
  @("sCs", "Fgf2", "Sdc4", "sCs", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 27.51115433333333, 82.533463, 0.4738307886277414, 0.4738307886277414, 61.80004835207445, 556.20043516867, 0.09625666637074914, 0.09625666637074916)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $data[$i]
  $rowNum = $i + 2
  for ($j = 0; $j -lt $r.Length; $j++) {
    $ws.Cells.Item($rowNum, $j + 1).Value = $r[$j]
  }
}